$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 290.6154
$ws.Cells.Item(33, 9).Value = 290.6154
$ws.Cells.Item(33, 11).Value = 290.6154
$ws.Cells.Item(33, 13).Value = -61.61540000000002

$ws.Cells.Item(53, 8).Value = 844.7
$ws.Cells.Item(53, 9).Value = 1164.7778
$ws.Cells.Item(53, 10).Value = 582.8182
$ws.Cells.Item(53, 11).Value = 1164.7778
$ws.Cells.Item(53, 12).Value = 582.8182
$ws.Cells.Item(53, 13).Value = -527.7778000000001
$ws.Cells.Item(53, 14).Value = -1856.8182

$ws.Cells.Item(116, 8).Value = 10399871
$ws.Cells.Item(116, 9).Value = 24906864
$ws.Cells.Item(116, 11).Value = 24906864
$ws.Cells.Item(116, 13).Value = -24903422

$ws.Cells.Item(132, 8).Value = 1916.825
$ws.Cells.Item(132, 9).Value = 1715.9744
$ws.Cells.Item(132, 11).Value = 5147.9232
$ws.Cells.Item(132, 13).Value = -2617.9232

$ws.Cells.Item(138, 8).Value = 2611.6736
$ws.Cells.Item(138, 9).Value = 2027.2174
$ws.Cells.Item(138, 10).Value = 3128.6924
$ws.Cells.Item(138, 11).Value = 6081.6522
$ws.Cells.Item(138, 12).Value = 9386.0772
$ws.Cells.Item(138, 13).Value = -941.6522000000004
$ws.Cells.Item(138, 14).Value = -19666.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 990000
$ws.Cells.Item(6, 9).Value = 990000
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 990000
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -989827
$ws.Cells.Item(6, 14).ClearContents() | Out-Null

$ws.Cells.Item(61, 8).Value = 29133.473
$ws.Cells.Item(61, 9).Value = 580.4400000000001
$ws.Cells.Item(61, 11).Value = 580.4400000000001
$ws.Cells.Item(61, 13).Value = -368.4400000000001

$ws.Cells.Item(132, 8).Value = 1171.1372
$ws.Cells.Item(132, 9).Value = 979.55554
$ws.Cells.Item(132, 11).Value = 2938.66662
$ws.Cells.Item(132, 13).Value = -408.66662

$ws.Cells.Item(136, 8).Value = 29133.473
$ws.Cells.Item(136, 9).Value = 580.4400000000001
$ws.Cells.Item(136, 11).Value = 1741.32
$ws.Cells.Item(136, 13).Value = 808.6799999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 3434.6667
$ws.Cells.Item(37, 9).Value = 3434.6667
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 3434.6667
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = -3297.6667
$ws.Cells.Item(37, 14).ClearContents() | Out-Null

$ws.Cells.Item(86, 8).Value = 7035.3335
$ws.Cells.Item(86, 9).Value = 4239
$ws.Cells.Item(86, 11).Value = 4239
$ws.Cells.Item(86, 13).Value = -3116

$ws.Cells.Item(89, 8).Value = 7035.3335
$ws.Cells.Item(89, 9).Value = 4239
$ws.Cells.Item(89, 11).Value = 21195
$ws.Cells.Item(89, 13).Value = -15579

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents() | Out-Null

$ws.Cells.Item(134, 8).Value = 1397.7142
$ws.Cells.Item(134, 9).Value = 844.9761999999999
$ws.Cells.Item(134, 11).Value = 2534.9286
$ws.Cells.Item(134, 13).Value = 0.07140000000026703

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2499.2727
$ws.Cells.Item(31, 9).Value = 2152.64
$ws.Cells.Item(31, 11).Value = 2152.64
$ws.Cells.Item(31, 13).Value = -1857.64

$ws.Cells.Item(34, 8).Value = 2499.2727
$ws.Cells.Item(34, 9).Value = 2152.64
$ws.Cells.Item(34, 11).Value = 2152.64
$ws.Cells.Item(34, 13).Value = -1950.64

$ws.Cells.Item(86, 8).Value = 7165.4375
$ws.Cells.Item(86, 9).Value = 5250
$ws.Cells.Item(86, 10).Value = 9080.875
$ws.Cells.Item(86, 11).Value = 5250
$ws.Cells.Item(86, 12).Value = 9080.875
$ws.Cells.Item(86, 13).Value = -4127
$ws.Cells.Item(86, 14).Value = -11326.875

$ws.Cells.Item(89, 8).Value = 7165.4375
$ws.Cells.Item(89, 9).Value = 5250
$ws.Cells.Item(89, 10).Value = 9080.875
$ws.Cells.Item(89, 11).Value = 26250
$ws.Cells.Item(89, 12).Value = 45404.375
$ws.Cells.Item(89, 13).Value = -20634
$ws.Cells.Item(89, 14).Value = -56636.375

$ws.Cells.Item(94, 8).Value = 23526.6
$ws.Cells.Item(94, 10).Value = 1649
$ws.Cells.Item(94, 12).Value = 1649
$ws.Cells.Item(94, 14).Value = -2551

$ws.Cells.Item(97, 8).Value = 44200
$ws.Cells.Item(97, 10).Value = 44200
$ws.Cells.Item(97, 12).Value = 44200
$ws.Cells.Item(97, 14).Value = -46182

$ws.Cells.Item(99, 8).Value = 12944848
$ws.Cells.Item(99, 9).Value = 15876246
$ws.Cells.Item(99, 11).Value = 15876246
$ws.Cells.Item(99, 13).Value = -15874748

$ws.Cells.Item(122, 8).Value = 3167.2
$ws.Cells.Item(122, 9).Value = 2632.5454
$ws.Cells.Item(122, 11).Value = 7897.6362
$ws.Cells.Item(122, 13).Value = -5447.6362

$ws.Cells.Item(126, 8).Value = 12944848
$ws.Cells.Item(126, 9).Value = 15876246
$ws.Cells.Item(126, 11).Value = 47628738
$ws.Cells.Item(126, 13).Value = -47626268

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 8583.691999999999
$ws.Cells.Item(114, 10).Value = 10936
$ws.Cells.Item(114, 12).Value = 32808
$ws.Cells.Item(114, 14).Value = -39316

$ws.Cells.Item(117, 8).Value = 2658.4285
$ws.Cells.Item(117, 10).Value = 2984.8333
$ws.Cells.Item(117, 12).Value = 8954.499899999999
$ws.Cells.Item(117, 14).Value = -15838.4999

$ws.Cells.Item(121, 8).Value = 2649.4
$ws.Cells.Item(121, 10).Value = 2959.8333
$ws.Cells.Item(121, 12).Value = 8879.499899999999
$ws.Cells.Item(121, 14).Value = -11499.4999

$ws.Cells.Item(131, 8).Value = 40999.8
$ws.Cells.Item(131, 9).Value = 63170.438
$ws.Cells.Item(131, 10).Value = 1585.3334
$ws.Cells.Item(131, 11).Value = 189511.314
$ws.Cells.Item(131, 12).Value = 4756.0002
$ws.Cells.Item(131, 13).Value = -184471.314
$ws.Cells.Item(131, 14).Value = -14836.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1352.9565
$ws.Cells.Item(102, 9).Value = 1352.9565
$ws.Cells.Item(102, 11).Value = 1352.9565
$ws.Cells.Item(102, 13).Value = 269.0435

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6182222
$ws.Cells.Item(40, 9).Value = 13485.2
$ws.Cells.Item(40, 11).Value = 13485.2
$ws.Cells.Item(40, 13).Value = -13349.2

$ws.Cells.Item(100, 8).Value = 11436.833
$ws.Cells.Item(100, 9).Value = 13091.333
$ws.Cells.Item(100, 11).Value = 13091.333
$ws.Cells.Item(100, 13).Value = -12550.333

$ws.Cells.Item(136, 8).Value = 1987.3529
$ws.Cells.Item(136, 9).Value = 1677.9166
$ws.Cells.Item(136, 11).Value = 5033.7498
$ws.Cells.Item(136, 13).Value = -2483.7498

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 18279.223
$ws.Cells.Item(45, 10).Value = 18279.223
$ws.Cells.Item(45, 12).Value = 18279.223
$ws.Cells.Item(45, 14).Value = -19261.223

$ws.Cells.Item(62, 8).Value = 8232.571
$ws.Cells.Item(62, 9).Value = 10225.385
$ws.Cells.Item(62, 11).Value = 10225.385
$ws.Cells.Item(62, 13).Value = -9601.385

$ws.Cells.Item(65, 8).Value = 8232.571
$ws.Cells.Item(65, 9).Value = 10225.385
$ws.Cells.Item(65, 11).Value = 51126.925
$ws.Cells.Item(65, 13).Value = -48006.925

$ws.Cells.Item(93, 8).Value = 73200.2
$ws.Cells.Item(93, 10).Value = 67500
$ws.Cells.Item(93, 12).Value = 67500
$ws.Cells.Item(93, 14).Value = -72492

$ws.Cells.Item(100, 8).Value = 5953068.5
$ws.Cells.Item(100, 9).Value = 6494188.5
$ws.Cells.Item(100, 11).Value = 12988377
$ws.Cells.Item(100, 13).Value = -12987836

$ws.Cells.Item(122, 8).Value = 2032.8206
$ws.Cells.Item(122, 10).Value = 2537.8462
$ws.Cells.Item(122, 12).Value = 7613.5386
$ws.Cells.Item(122, 14).Value = -12513.5386

$ws.Cells.Item(136, 8).Value = 971.2258
$ws.Cells.Item(136, 9).Value = 895.26666
$ws.Cells.Item(136, 10).Value = 3250
$ws.Cells.Item(136, 11).Value = 2685.79998
$ws.Cells.Item(136, 12).Value = 9750
$ws.Cells.Item(136, 13).Value = -135.7999799999998
$ws.Cells.Item(136, 14).Value = -14850
